$d = $word.ActiveDocument

$replacements = @(
    @("48÷9=5, 3", "75÷9=8, 3"),
    @("44÷5=8, 4", "45÷7=6, 3"),
    @("98÷5=19, 3", "67÷4=16, 3"),
    @("42÷2=21, 0", "65÷2=32, 1"),
    @("70÷6=11, 4", "13÷7=1, 6"),
    @("34÷3=11, 1", "33÷2=16, 1"),
    @("26÷3=8, 2", "33÷4=8, 1"),
    @("44÷7=6, 2", "42÷4=10, 2"),
    @("53÷8=6, 5", "41÷2=20, 1"),
    @("43÷8=5, 3", "61÷6=10, 1"),
    @("15÷4=3, 3", "75÷8=9, 3"),
    @("79÷4=19, 3", "25÷5=5, 0"),
    @("73÷8=9, 1", "29÷6=4, 5"),
    @("37÷3=12, 1", "23÷5=4, 3"),
    @("34÷8=4, 2", "84÷3=28, 0"),
    @("66÷2=33, 0", "10÷2=5, 0"),
    @("75÷4=18, 3", "37÷4=9, 1"),
    @("65÷3=21, 2", "56÷4=14, 0"),
    @("95÷3=31, 2", "97÷5=19, 2"),
    @("52÷8=6, 4", "51÷4=12, 3"),
    @("60÷5=12, 0", "83÷8=10, 3"),
    @("92÷3=30, 2", "74÷4=18, 2"),
    @("20÷8=2, 4", "27÷6=4, 3"),
    @("61÷5=12, 1", "86÷6=14, 2"),
    @("48÷7=6, 6", "70÷5=14, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
